# Update "想去人数" (number of people interested) values on the
# "展览" and "全部类型" sheets to match the newly published data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 2313
    $ws.Range("F3").Value = 1767
    $ws.Range("F6").Value = 957
    $ws.Range("F8").Value = 5877
    $ws.Range("F9").Value = 92
}
